$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing May-11 total_venda value
$ws.Range("B12").Value = 9776.9

# Insert a new row for May day 12, shifting everything below down by one
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 31185.44
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 2025
$ws.Range("E13").Value = "05/2025"
